$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @{
    17 = @('624042702533402', 'true', '0', '0', '0', '237669800046', 'false', 'None', 'None', 'false', 'None', 'None', 'None', 'None', 'None', 'None', 'UNKNOWN SUBSCRIBER', 'UNKNOWN SUBSCRIBER', 'UNKNOWN SUBSCRIBER', 'UNKNOWN SUBSCRIBER', 'isActiveIMSI:Subscriber is not attached. Try to restart phone;')
    18 = @('624042732168871', 'true', '0', '0', '0', '237664022676', 'true', '3566780923618278', 'None', 'true', '237660002052', '10.124.148.4', 'None', '160664022676', '160664022676', '160664022676', 'UNKNOWN SUBSCRIBER', 'UNKNOWN SUBSCRIBER', 'UNKNOWN SUBSCRIBER', 'KNOWN SUBSCRIBER', 'DOMS02:KNOWN SUBSCRIBER;cfb:cfb is defined to 160664022676;cfnrc:cfnrc is defined to 160664022676;cfnry:cfnry is defined to 160664022676;')
    19 = @('624042732168871', 'true', '1', '0', '0', '237664022676', 'true', '3566780923618278', 'None', 'true', '237660002052', '10.124.148.4', 'None', '160664022676', '160664022676', '160664022676', 'UNKNOWN SUBSCRIBER', 'UNKNOWN SUBSCRIBER', 'UNKNOWN SUBSCRIBER', 'KNOWN SUBSCRIBER', 'odboc:Barring oc in HLR;DOMS02:KNOWN SUBSCRIBER;cfb:cfb is defined to 160664022676;cfnrc:cfnrc is defined to 160664022676;cfnry:cfnry is defined to 160664022676;')
    20 = @('624042732168871', 'true', '3', '0', '0', '237664022676', 'true', '3566780923618278', 'None', 'true', '237660002052', '10.124.148.4', 'None', '160664022676', '160664022676', '160664022676', 'UNKNOWN SUBSCRIBER', 'UNKNOWN SUBSCRIBER', 'UNKNOWN SUBSCRIBER', 'KNOWN SUBSCRIBER', 'odboc:Barring oc in HLR;DOMS02:KNOWN SUBSCRIBER;cfb:cfb is defined to 160664022676;cfnrc:cfnrc is defined to 160664022676;cfnry:cfnry is defined to 160664022676;')
    21 = @('624042732168871', 'true', '0', '0', '0', '237664022676', 'true', '3566780923618278', 'None', 'true', '237660002052', '10.124.148.4', 'None', '160664022676', '16066022676', '160664022676', 'UNKNOWN SUBSCRIBER', 'UNKNOWN SUBSCRIBER', 'UNKNOWN SUBSCRIBER', 'KNOWN SUBSCRIBER', 'DOMS02:KNOWN SUBSCRIBER;cfb:cfb is defined to 160664022676;cfnrc:cfnrc is defined to 16066022676;cfnry:cfnry is defined to 160664022676;')
    22 = @('624042732168871', 'true', '0', '0', '0', '237664022676', 'true', '3566780923618278', 'None', 'true', '237660002052', '10.124.148.4', 'None', '160664022676', '16066022676', '160664022676', 'UNKNOWN SUBSCRIBER', 'UNKNOWN SUBSCRIBER', 'UNKNOWN SUBSCRIBER', 'KNOWN SUBSCRIBER', 'DOMS02:KNOWN SUBSCRIBER;cfb:cfb is defined to 160664022676;cfnrc:cfnrc is defined to 16066022676;cfnry:cfnry is defined to 160664022676;')
    23 = @('624042732168871', 'true', '0', '0', '0', '237664022676', 'true', '3566780923618278', 'None', 'false', 'None', 'None', 'None', '160664022676', '160664022676', '160664022676', 'UNKNOWN SUBSCRIBER', 'UNKNOWN SUBSCRIBER', 'UNKNOWN SUBSCRIBER', 'UNKNOWN SUBSCRIBER', 'no_exist:Subscriber exists in any MSS. Restart phone;cfb:cfb is defined to 160664022676;cfnrc:cfnrc is defined to 160664022676;cfnry:cfnry is defined to 160664022676;')
    24 = @('624042700253500', 'true', '0', '0', '0', '237661341827', 'true', '8655500200407700', 'None', 'true', '237660001052', '10.124.208.81', 'None', 'None', 'None', 'None', 'UNKNOWN SUBSCRIBER', 'UNKNOWN SUBSCRIBER', 'UNKNOWN SUBSCRIBER', 'UNKNOWN SUBSCRIBER', 'no_exist:Subscriber exists in any MSS. Restart phone;result:Everything is ok in HLR;')
    25 = @('624042702533402', 'true', '0', '0', '0', '237669800046', 'false', 'None', 'None', 'false', 'None', 'None', 'None', 'None', 'None', 'None', 'UNKNOWN SUBSCRIBER', 'UNKNOWN SUBSCRIBER', 'UNKNOWN SUBSCRIBER', 'UNKNOWN SUBSCRIBER', 'isActiveIMSI:Subscriber is not attached. Try to restart phone;')
}

foreach ($r in $newRows.Keys | Sort-Object {[int]$_}) {
    $vals = $newRows[$r]
    for ($c = 1; $c -le $vals.Length; $c++) {
        $cell = $ws.Cells.Item([int]$r, $c)
        $cell.NumberFormat = "@"
        $cell.Value = "'" + $vals[$c - 1]
    }
}
